$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1078.1364
$ws.Range("I17").Value = 651.5
$ws.Range("J17").Value = 1238.125
$ws.Range("K17").Value = 1954.5
$ws.Range("L17").Value = 3714.375
$ws.Range("M17").Value = -1786.5
$ws.Range("N17").Value = -4050.375
$ws.Range("H28").Value = 192
$ws.Range("I28").Value = 128.10715
$ws.Range("J28").Value = 549.8
$ws.Range("K28").Value = 128.10715
$ws.Range("L28").Value = 549.8
$ws.Range("M28").Value = 356.89285
$ws.Range("N28").Value = -1519.8
$ws.Range("H33").Value = 983.1429000000001
$ws.Range("I33").Value = 793.375
$ws.Range("J33").Value = 1590.4
$ws.Range("K33").Value = 793.375
$ws.Range("L33").Value = 1590.4
$ws.Range("M33").Value = -564.375
$ws.Range("N33").Value = -2048.4
$ws.Range("H111").Value = 738.1667
$ws.Range("I111").Value = 485.8
$ws.Range("J111").Value = 2000
$ws.Range("K111").Value = 1457.4
$ws.Range("L111").Value = 6000
$ws.Range("M111").Value = 1609.6
$ws.Range("N111").Value = -12134
$ws.Range("H112").Value = 1842.1428
$ws.Range("I112").Value = 500
$ws.Range("J112").Value = 1945.3846
$ws.Range("K112").Value = 1500
$ws.Range("L112").Value = 5836.1538
$ws.Range("M112").Value = -392
$ws.Range("N112").Value = -8052.1538
$ws.Range("H113").Value = 3496.394
$ws.Range("I113").Value = 3156.818
$ws.Range("J113").Value = 3666.182
$ws.Range("K113").Value = 3156.818
$ws.Range("L113").Value = 3666.182
$ws.Range("M113").Value = 97.18199999999979
$ws.Range("N113").Value = -10174.182
$ws.Range("H135").Value = 659.8125
$ws.Range("I135").Value = 653.8182
$ws.Range("J135").Value = 673
$ws.Range("K135").Value = 5884.3638
$ws.Range("L135").Value = 6057
$ws.Range("M135").Value = -3349.3638
$ws.Range("N135").Value = -11127
$ws.Range("H138").Value = 2632.096
$ws.Range("I138").Value = 1922.7142
$ws.Range("J138").Value = 2918.577
$ws.Range("K138").Value = 5768.142599999999
$ws.Range("L138").Value = 8755.731
$ws.Range("M138").Value = -628.1425999999992
$ws.Range("N138").Value = -19035.731
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1840
$ws.Range("I45").Value = 1723.3334
$ws.Range("J45").Value = 1980
$ws.Range("K45").Value = 1723.3334
$ws.Range("L45").Value = 1980
$ws.Range("M45").Value = -1346.3334
$ws.Range("N45").Value = -2734
$ws.Range("H97").Value = 1320
$ws.Range("I97").Value = 1268
$ws.Range("J97").Value = 1476
$ws.Range("K97").Value = 1268
$ws.Range("L97").Value = 1476
$ws.Range("M97").Value = -772
$ws.Range("N97").Value = -2468
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5115.567
$ws.Range("I86").Value = 5113.8
$ws.Range("J86").Value = 5117.3335
$ws.Range("K86").Value = 5113.8
$ws.Range("L86").Value = 5117.3335
$ws.Range("M86").Value = -3990.8
$ws.Range("N86").Value = -7363.3335
$ws.Range("H89").Value = 5115.567
$ws.Range("I89").Value = 5113.8
$ws.Range("J89").Value = 5117.3335
$ws.Range("K89").Value = 25569
$ws.Range("L89").Value = 25586.6675
$ws.Range("M89").Value = -19953
$ws.Range("N89").Value = -36818.6675
$ws.Range("H99").Value = 6289
$ws.Range("I99").Value = 6289
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 6289
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = ""
$ws.Range("N99").Value = -4791
$ws.Range("H105").Value = 2710.22
$ws.Range("I105").Value = 2443.2432
$ws.Range("J105").Value = 3470.077
$ws.Range("K105").Value = 2443.2432
$ws.Range("L105").Value = 3470.077
$ws.Range("M105").Value = -696.2431999999999
$ws.Range("N105").Value = -6964.077
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 167466.67
$ws.Range("I16").Value = 167466.67
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 167466.67
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = ""
$ws.Range("N16").Value = -167179.67
$ws.Range("H31").Value = 3083.426
$ws.Range("I31").Value = 2399.0303
$ws.Range("J31").Value = 4158.905
$ws.Range("K31").Value = 2399.0303
$ws.Range("L31").Value = 4158.905
$ws.Range("M31").Value = -2104.0303
$ws.Range("N31").Value = -4748.905
$ws.Range("H34").Value = 3083.426
$ws.Range("I34").Value = 2399.0303
$ws.Range("J34").Value = 4158.905
$ws.Range("K34").Value = 2399.0303
$ws.Range("L34").Value = 4158.905
$ws.Range("M34").Value = -2197.0303
$ws.Range("N34").Value = -4562.905
$ws.Range("H58").Value = 2533.4119
$ws.Range("I58").Value = 1894.3
$ws.Range("J58").Value = 3446.4285
$ws.Range("K58").Value = 1894.3
$ws.Range("L58").Value = 3446.4285
$ws.Range("M58").Value = -1691.3
$ws.Range("N58").Value = -3852.4285
$ws.Range("H62").Value = 4993.778
$ws.Range("I62").Value = 3978
$ws.Range("J62").Value = 6263.5
$ws.Range("K62").Value = 3978
$ws.Range("L62").Value = 6263.5
$ws.Range("M62").Value = -3354
$ws.Range("N62").Value = -7511.5
$ws.Range("H65").Value = 4993.778
$ws.Range("I65").Value = 3978
$ws.Range("J65").Value = 6263.5
$ws.Range("K65").Value = 19890
$ws.Range("L65").Value = 31317.5
$ws.Range("M65").Value = -16770
$ws.Range("N65").Value = -37557.5
$ws.Range("H105").Value = 625.6094000000001
$ws.Range("I105").Value = 611.6070999999999
$ws.Range("J105").Value = 723.625
$ws.Range("K105").Value = 611.6070999999999
$ws.Range("L105").Value = 723.625
$ws.Range("M105").Value = 1135.3929
$ws.Range("N105").Value = -4217.625
$ws.Range("H113").Value = 167466.67
$ws.Range("I113").Value = 167466.67
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 167466.67
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = ""
$ws.Range("N113").Value = -165296.67
$ws.Range("H136").Value = 2533.4119
$ws.Range("I136").Value = 1894.3
$ws.Range("J136").Value = 3446.4285
$ws.Range("K136").Value = 5682.9
$ws.Range("L136").Value = 10339.2855
$ws.Range("M136").Value = -3132.9
$ws.Range("N136").Value = -15439.2855
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 951.28
$ws.Range("I97").Value = 743.38464
$ws.Range("J97").Value = 1176.5
$ws.Range("K97").Value = 743.38464
$ws.Range("L97").Value = 1176.5
$ws.Range("M97").Value = -247.38464
$ws.Range("N97").Value = -2168.5
$ws.Range("H135").Value = 48800
$ws.Range("J135").Value = 50000
$ws.Range("L135").Value = 50000
$ws.Range("N135").Value = -60140
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2713.3
$ws.Range("I7").Value = 2735
$ws.Range("J7").Value = 2691.6
$ws.Range("K7").Value = 2735
$ws.Range("L7").Value = 2691.6
$ws.Range("M7").Value = -2623
$ws.Range("N7").Value = -2915.6
$ws.Range("H126").Value = 2713.3
$ws.Range("I126").Value = 2735
$ws.Range("J126").Value = 2691.6
$ws.Range("K126").Value = 8205
$ws.Range("L126").Value = 8074.799999999999
$ws.Range("M126").Value = -5735
$ws.Range("N126").Value = -13014.8
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 7501.5
$ws.Range("I96").Value = 3001.5
$ws.Range("J96").Value = 12001.5
$ws.Range("K96").Value = 3001.5
$ws.Range("L96").Value = 12001.5
$ws.Range("M96").Value = -1628.5
$ws.Range("N96").Value = -14747.5
$ws.Range("H100").Value = 3378.3914
$ws.Range("I100").Value = 1325.6
$ws.Range("J100").Value = 17063.666
$ws.Range("K100").Value = 2651.2
$ws.Range("L100").Value = 34127.332
$ws.Range("M100").Value = -2110.2
$ws.Range("N100").Value = -35209.332
$ws.Range("H107").Value = 660.7857
$ws.Range("I107").Value = 794.5714
$ws.Range("J107").Value = 527
$ws.Range("K107").Value = 2383.7142
$ws.Range("L107").Value = 1581
$ws.Range("M107").Value = -463.7142000000003
$ws.Range("N107").Value = -5421

Write-Output "applied changes"